$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column headers
$ws.Range("E1").Value = "tiempo_d-r"
$ws.Range("F1").Value = "presion_d-r(mTorr)"

# New data for columns E (tiempo_d-r) and F (presion_d-r(mTorr))
$dataE = @(58, 60, 65, 70, 75, 80, 85, 90, 95, 100, 105, 110, 120, 135, 153)
$dataF = @(260, 240, 220, 260, 220, 170, 160, 150, 140, 140, 140, 130, 130, 120, 120)

for ($i = 0; $i -lt $dataE.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $dataE[$i]
    $ws.Cells.Item($row, 6).Value = $dataF[$i]
}

$ws.Range("E15").Select()
